$data = @(
    ,@('510880','华泰柏瑞上证红利ETF','181.00','97.22','3.19','5.7739',3)
    ,@('519918','华夏兴和混合','53.07','82.18','4.82','2.5580',6)
    ,@('160311','华夏蓝筹混合(LOF)','30.12','87.73','5.48','1.6506',5)
    ,@('100032','富国中证红利指数增强','56.72','94.48','1.87','1.0607',7)
    ,@('090010','大成中证红利指数A','34.51','93.73','1.73','0.5970',5)
    ,@('515180','易方达中证红利ETF','16.55','99.58','1.84','0.3045',5)
    ,@('501029','华宝标普中国A股红利机会指数（LOF）A','13.19','94.39','1.35','0.1781',9)
    ,@('515080','招商中证红利ETF','9.06','99.25','1.83','0.1658',5)
    ,@('007801','大成中证红利指数C','3.87','93.73','1.73','0.0670',5)
    ,@('005562','创金合信中证红利低波动指数C','1.85','94.46','3.45','0.0638',3)
    ,@('005576','华泰柏瑞新金融地产灵活配置混合','0.79','94.50','6.43','0.0508',7)
    ,@('512890','华泰柏瑞中证红利低波动ETF','1.36','99.24','3.62','0.0492',3)
    ,@('515890','博时中证红利ETF','2.59','98.55','1.82','0.0471',5)
    ,@('512040','富国中证价值ETF','3.44','99.55','1.29','0.0444',7)
    ,@('005561','创金合信中证红利低波动指数A','1.22','94.46','3.45','0.0421',3)
    ,@('001244','华泰柏瑞量化智慧灵活配置混合A','3.50','91.02','1.06','0.0371',2)
    ,@('007518','东方阿尔法优选混合A','2.03','72.64','1.60','0.0325',6)
    ,@('009726','招商中证500等权重指数增强A','1.87','91.11','1.56','0.0292',2)
    ,@('161907','万家中证红利指数(LOF)','1.34','94.87','1.74','0.0233',5)
    ,@('900027','中信证券信远一年持有期混合型集合资产管理计划A','0.71','75.94','3.01','0.0214',8)
    ,@('006123','中融高股息精选混合A','0.58','92.22','3.65','0.0212',4)
    ,@('006652','富国金融地产行业混合A','0.45','90.56','4.52','0.0203',9)
    ,@('007519','东方阿尔法优选混合C','0.82','72.64','1.60','0.0131',6)
    ,@('009727','招商中证500等权重指数增强C','0.69','91.11','1.56','0.0108',2)
    ,@('006104','华泰柏瑞量化智慧灵活配置混合C','0.84','91.02','1.06','0.0089',2)
    ,@('006124','中融高股息精选混合C','0.23','92.22','3.65','0.0084',4)
    ,@('001614','东方区域发展混合','0.06','94.77','5.41','0.0032',7)
    ,@('519117','浦银安盛基本面400指数','0.24','92.63','0.79','0.0019',2)
    ,@('162907','泰信中证锐联基本面400指数（LOF）','0.23','94.61','0.77','0.0018',4)
    ,@('011124','富国金融地产行业混合C','0.03','90.56','4.52','0.0014',9)
    ,@('001273','民生加银新动力灵活配置混合A','0.04','68.44','2.12','0.0008',5)
    ,@('001274','民生加银新动力灵活配置混合D','0.04','68.44','2.12','0.0008',5)
    ,@('900087','中信证券信远一年持有期混合型集合资产管理计划C','0.02','75.94','3.01','0.0006',8)
    ,@('900077','中信证券信远一年持有期混合型集合资产管理计划B','0.01','75.94','3.01','0.0003',8)
)
# ============================================================================
# Edit: add a new "2022-Q1" worksheet (fund-holding detail) before the
# "总计" (totals) sheet, and append a corresponding summary row to "总计".
# ============================================================================

$wb = $excel.ActiveWorkbook

# ---- locate reference sheets -------------------------------------------
$templateSheet = $wb.Worksheets.Item("2021-Q4")
$insertBefore = $wb.Worksheets.Item("总计")

# ---- create the new sheet, positioned immediately before "总计" ---------
$newSheet = $wb.Worksheets.Add($insertBefore)
$newSheet.Name = "2022-Q1"

# NOTE: the worksheet handle obtained before the Add() call tracks a sheet
# *position*, so after inserting a new sheet at that position the old
# handle now refers to the newly inserted sheet instead of "总计". Re-fetch
# the "总计" worksheet by name now that the layout is final.
$totalSheet = $wb.Worksheets.Item("总计")

# ---- header row (B1:H1) --------------------------------------------------
$templateSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$newSheet.Cells.Item(1,2).Value = "基金代码"
$newSheet.Cells.Item(1,3).Value = "基金名称"
$newSheet.Cells.Item(1,4).Value = "基金规模"
$newSheet.Cells.Item(1,5).Value = "股票总仓位"
$newSheet.Cells.Item(1,6).Value = "仓位占比"
$newSheet.Cells.Item(1,7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1,8).Value = "仓位排名"


# ---- data rows (A2:H35) --------------------------------------------------
$rowCount = $data.Count

# Style column A (index numbers) the same way the template sheet does.
$templateSheet.Range("A2").Copy()
$newSheet.Range("A2:A" + (1 + $rowCount)).PasteSpecial(-4122)

for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $i + 2
    $row = $data[$i]

    $newSheet.Cells.Item($r,1).Value = $i

    # Force text for the string-valued columns B..G so that values such as
    # "181.00" / "0.0670" are preserved exactly (not coerced to numbers).
    $newSheet.Cells.Item($r,2).Value = "'" + $row[0]
    $newSheet.Cells.Item($r,3).Value = "'" + $row[1]
    $newSheet.Cells.Item($r,4).Value = "'" + $row[2]
    $newSheet.Cells.Item($r,5).Value = "'" + $row[3]
    $newSheet.Cells.Item($r,6).Value = "'" + $row[4]
    $newSheet.Cells.Item($r,7).Value = "'" + $row[5]

    # Numeric ranking column.
    $newSheet.Cells.Item($r,8).Value = $row[6]
}

# Drop the "quote prefix" formatting introduced above so the cells end up
# with no explicit style, matching plain text cells.
$newSheet.Range("B2:G" + (1 + $rowCount)).ClearFormats()


# ---- update the "总计" (totals) sheet ------------------------------------
# Insert a new row right after the header for the 2022-Q1 summary, shifting
# the previously existing quarters down by one.
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()

# Re-apply the index-column style (copied from the row pushed down to A3).
$totalSheet.Cells.Item(3,1).Copy()
$totalSheet.Cells.Item(2,1).PasteSpecial(-4122)

$totalSheet.Cells.Item(2,1).Value = 0
$totalSheet.Cells.Item(2,2).Value = "2022-Q1"
$totalSheet.Cells.Item(2,3).Value = 34
$totalSheet.Cells.Item(2,4).Value = 12.89

# Renumber the index column for the rows that shifted down.
$lastRow = $totalSheet.UsedRange.Rows.Count
for ($r = 3; $r -le $lastRow; $r++) {
    $totalSheet.Cells.Item($r,1).Value = $r - 2
}
